# Applies the numeric updates described by the commit diff to the
# "Jenova_Profits" aggregated price/profit workbook (8 sheets: ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR). Each touched row represents one Leve,
# and columns H..N hold price/profit figures that were refreshed by the
# scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 763.75
$ws.Cells.Item(33, 9).Value = 814.8182
$ws.Cells.Item(33, 11).Value = 814.8182
$ws.Cells.Item(33, 13).Value = -585.8182

$ws.Cells.Item(45, 8).Value = 3482.6667
$ws.Cells.Item(45, 10).Value = 3793
$ws.Cells.Item(45, 12).Value = 11379
$ws.Cells.Item(45, 14).Value = -11763

$ws.Cells.Item(86, 8).Value = 6582321
$ws.Cells.Item(86, 9).Value = 4218.4
$ws.Cells.Item(86, 11).Value = 4218.4
$ws.Cells.Item(86, 13).Value = -3095.4

$ws.Cells.Item(89, 8).Value = 6582321
$ws.Cells.Item(89, 9).Value = 4218.4
$ws.Cells.Item(89, 11).Value = 21092
$ws.Cells.Item(89, 13).Value = -15476

$ws.Cells.Item(106, 8).Value = 2908.1177
$ws.Cells.Item(106, 9).Value = 3643.5
$ws.Cells.Item(106, 10).Value = 2254.4443
$ws.Cells.Item(106, 11).Value = 3643.5
$ws.Cells.Item(106, 12).Value = 2254.4443
$ws.Cells.Item(106, 13).Value = -3012.5
$ws.Cells.Item(106, 14).Value = -3516.4443

$ws.Cells.Item(107, 8).Value = 37602.223
$ws.Cells.Item(107, 9).Value = 48175.285
$ws.Cells.Item(107, 11).Value = 48175.285
$ws.Cells.Item(107, 13).Value = -46255.285

$ws.Cells.Item(108, 8).Value = 70000
$ws.Cells.Item(108, 10).Value = 70000
$ws.Cells.Item(108, 12).Value = 70000
$ws.Cells.Item(108, 14).Value = -77680

$ws.Cells.Item(125, 8).Value = 3990.625
$ws.Cells.Item(125, 9).Value = 3197.8
$ws.Cells.Item(125, 11).Value = 28780.2
$ws.Cells.Item(125, 13).Value = -26320.2

$ws.Cells.Item(127, 8).Value = 34743.75
$ws.Cells.Item(127, 9).Value = 34743.75
$ws.Cells.Item(127, 11).Value = 104231.25
$ws.Cells.Item(127, 13).Value = -99271.25

$ws.Cells.Item(129, 8).Value = 15718.8125
$ws.Cells.Item(129, 9).Value = 38003.668
$ws.Cells.Item(129, 10).Value = 10576.154
$ws.Cells.Item(129, 11).Value = 114011.004
$ws.Cells.Item(129, 12).Value = 31728.462
$ws.Cells.Item(129, 13).Value = -109011.004
$ws.Cells.Item(129, 14).Value = -41728.462

$ws.Cells.Item(138, 8).Value = 5880.12
$ws.Cells.Item(138, 9).Value = 3320.9333
$ws.Cells.Item(138, 10).Value = 6976.914
$ws.Cells.Item(138, 11).Value = 9962.7999
$ws.Cells.Item(138, 12).Value = 20930.742
$ws.Cells.Item(138, 13).Value = -4822.7999
$ws.Cells.Item(138, 14).Value = -31210.742

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 24000
$ws.Cells.Item(43, 10).Value = 24000
$ws.Cells.Item(43, 12).Value = 24000
$ws.Cells.Item(43, 14).Value = -24626

$ws.Cells.Item(45, 8).Value = 3397.7727
$ws.Cells.Item(45, 10).Value = 3276.4
$ws.Cells.Item(45, 12).Value = 3276.4
$ws.Cells.Item(45, 14).Value = -4030.4

$ws.Cells.Item(61, 8).Value = 3998.261
$ws.Cells.Item(61, 9).Value = 2531
$ws.Cells.Item(61, 11).Value = 2531
$ws.Cells.Item(61, 13).Value = -2319

$ws.Cells.Item(103, 8).Value = 35362
$ws.Cells.Item(103, 10).Value = 35362
$ws.Cells.Item(103, 12).Value = 35362
$ws.Cells.Item(103, 14).Value = -37706

$ws.Cells.Item(136, 8).Value = 3998.261
$ws.Cells.Item(136, 9).Value = 2531
$ws.Cells.Item(136, 11).Value = 7593
$ws.Cells.Item(136, 13).Value = -5043

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2158.8572
$ws.Cells.Item(20, 9).Value = 1844.0588
$ws.Cells.Item(20, 11).Value = 1844.0588
$ws.Cells.Item(20, 13).Value = -1597.0588

$ws.Cells.Item(94, 8).Value = 283.25
$ws.Cells.Item(94, 9).Value = 224.10527
$ws.Cells.Item(94, 10).Value = 408.1111
$ws.Cells.Item(94, 11).Value = 224.10527
$ws.Cells.Item(94, 12).Value = 408.1111
$ws.Cells.Item(94, 13).Value = 226.89473
$ws.Cells.Item(94, 14).Value = -1310.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 241238.12
$ws.Cells.Item(31, 9).Value = 835482.25
$ws.Cells.Item(31, 11).Value = 835482.25
$ws.Cells.Item(31, 13).Value = -835187.25

$ws.Cells.Item(34, 8).Value = 241238.12
$ws.Cells.Item(34, 9).Value = 835482.25
$ws.Cells.Item(34, 11).Value = 835482.25
$ws.Cells.Item(34, 13).Value = -835280.25

$ws.Cells.Item(86, 8).Value = 7902.222
$ws.Cells.Item(86, 10).Value = 8596.5
$ws.Cells.Item(86, 12).Value = 8596.5
$ws.Cells.Item(86, 14).Value = -10842.5

$ws.Cells.Item(89, 8).Value = 7902.222
$ws.Cells.Item(89, 10).Value = 8596.5
$ws.Cells.Item(89, 12).Value = 42982.5
$ws.Cells.Item(89, 14).Value = -54214.5

$ws.Cells.Item(99, 8).Value = 4974.12
$ws.Cells.Item(99, 9).Value = 3667.2856
$ws.Cells.Item(99, 11).Value = 3667.2856
$ws.Cells.Item(99, 13).Value = -2169.2856

$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()  # N110: remove (was -78178.836)

$ws.Cells.Item(122, 8).Value = 2462.96
$ws.Cells.Item(122, 9).Value = 1767.8823
$ws.Cells.Item(122, 10).Value = 3940
$ws.Cells.Item(122, 11).Value = 5303.6469
$ws.Cells.Item(122, 12).Value = 11820
$ws.Cells.Item(122, 13).Value = -2853.6469
$ws.Cells.Item(122, 14).Value = -16720

$ws.Cells.Item(126, 8).Value = 4974.12
$ws.Cells.Item(126, 9).Value = 3667.2856
$ws.Cells.Item(126, 11).Value = 11001.8568
$ws.Cells.Item(126, 13).Value = -8531.856800000001

$ws.Cells.Item(134, 8).Value = 3352.74
$ws.Cells.Item(134, 9).Value = 2330.9697
$ws.Cells.Item(134, 11).Value = 6992.909100000001
$ws.Cells.Item(134, 13).Value = -4457.909100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 1500
$ws.Cells.Item(7, 9).Value = 1000
$ws.Cells.Item(7, 10).Value = 2000
$ws.Cells.Item(7, 11).Value = 3000
$ws.Cells.Item(7, 12).Value = 6000
$ws.Cells.Item(7, 13).Value = -2888
$ws.Cells.Item(7, 14).Value = -6224

$ws.Cells.Item(44, 8).Value = 65.5
$ws.Cells.Item(44, 10).Value = 75.25
$ws.Cells.Item(44, 12).Value = 225.75
$ws.Cells.Item(44, 14).Value = -1021.75

$ws.Cells.Item(50, 8).Value = 11671
$ws.Cells.Item(50, 9).Value = 207.8
$ws.Cells.Item(50, 10).Value = 26000
$ws.Cells.Item(50, 11).Value = 623.4000000000001
$ws.Cells.Item(50, 12).Value = 78000
$ws.Cells.Item(50, 13).Value = -142.4000000000001
$ws.Cells.Item(50, 14).Value = -78962

$ws.Cells.Item(53, 8).Value = 11671
$ws.Cells.Item(53, 9).Value = 207.8
$ws.Cells.Item(53, 10).Value = 26000
$ws.Cells.Item(53, 11).Value = 623.4000000000001
$ws.Cells.Item(53, 12).Value = 78000
$ws.Cells.Item(53, 13).Value = -142.4000000000001
$ws.Cells.Item(53, 14).Value = -78962

$ws.Cells.Item(99, 8).Value = 3907.6667
$ws.Cells.Item(99, 9).Value = 400
$ws.Cells.Item(99, 10).Value = 4609.2
$ws.Cells.Item(99, 11).Value = 1200
$ws.Cells.Item(99, 12).Value = 13827.6
$ws.Cells.Item(99, 13).Value = 1046
$ws.Cells.Item(99, 14).Value = -18319.6

$ws.Cells.Item(108, 8).Value = 4680.091
$ws.Cells.Item(108, 9).Value = 4948.1
$ws.Cells.Item(108, 11).Value = 14844.3
$ws.Cells.Item(108, 13).Value = -11964.3

$ws.Cells.Item(119, 8).Value = 1724.75
$ws.Cells.Item(119, 9).Value = 966.3333
$ws.Cells.Item(119, 11).Value = 2898.9999
$ws.Cells.Item(119, 13).Value = 1939.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 14).ClearContents()  # N33: remove (was -14754)

$ws.Cells.Item(36, 8).Value = 16670439
$ws.Cells.Item(36, 9).Value = 25002158
$ws.Cells.Item(36, 10).Value = 7000
$ws.Cells.Item(36, 11).Value = 25002158
$ws.Cells.Item(36, 12).Value = 7000
$ws.Cells.Item(36, 13).Value = -25001673
$ws.Cells.Item(36, 14).Value = -7970

$ws.Cells.Item(132, 8).Value = 252545.03
$ws.Cells.Item(132, 9).Value = 314111.06
$ws.Cells.Item(132, 11).Value = 942333.1799999999
$ws.Cells.Item(132, 13).Value = -939803.1799999999

$ws.Cells.Item(136, 8).Value = 78330.664
$ws.Cells.Item(136, 10).Value = 78330.664
$ws.Cells.Item(136, 12).Value = 234991.992
$ws.Cells.Item(136, 14).Value = -240091.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5773.5557
$ws.Cells.Item(40, 9).Value = 4660.3335
$ws.Cells.Item(40, 10).Value = 8000
$ws.Cells.Item(40, 11).Value = 4660.3335
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = -4524.3335
$ws.Cells.Item(40, 14).Value = -8272

$ws.Cells.Item(105, 8).Value = 7615
$ws.Cells.Item(105, 10).Value = 7615
$ws.Cells.Item(105, 12).Value = 7615
$ws.Cells.Item(105, 14).Value = -14603

$ws.Cells.Item(136, 8).Value = 3363.1943
$ws.Cells.Item(136, 9).Value = 2414.1177
$ws.Cells.Item(136, 10).Value = 4212.3687
$ws.Cells.Item(136, 11).Value = 7242.353099999999
$ws.Cells.Item(136, 12).Value = 12637.1061
$ws.Cells.Item(136, 13).Value = -4692.353099999999
$ws.Cells.Item(136, 14).Value = -17737.1061

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1264.5555
$ws.Cells.Item(113, 10).Value = 2015.5
$ws.Cells.Item(113, 12).Value = 6046.5
$ws.Cells.Item(113, 14).Value = -10386.5

$ws.Cells.Item(126, 8).Value = 3323
$ws.Cells.Item(126, 9).Value = 2936.1428
$ws.Cells.Item(126, 11).Value = 8808.428400000001
$ws.Cells.Item(126, 13).Value = -6338.428400000001
